$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark that currently sits in the Pre-v1.0
#    paragraph ("... to find and read/write to addresses directly.").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. "3.0-beta1" -> "3.0-beta2", split across two runs ("-beta" and "2") to
#    match how Word represents the edited text.
$betaRange = $d.Content
$found = $betaRange.Find.Execute("-beta1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $betaRange.Text = "-beta"
    $afterBeta = $d.Range($betaRange.End, $betaRange.End)
    $afterBeta.InsertAfter("2")
}

# 3. "Check the alpha release notes on " -> "Check the release notes on "
$d.Content.Find.Execute("Check the alpha release notes on ", $true, $false, $false, $false, $false, $true, 1, $false, "Check the release notes on ", 2) | Out-Null

# 4. Insert two new bullet paragraphs before "Fix other cars not moving when
#    clutch pressed".
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Fix other cars not moving when clutch pressed*") {
        $p.Range.InsertBefore("Fix vehicle change detection`r")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Fix other cars not moving when clutch pressed*") {
        $p.Range.InsertBefore("Fix specific first gear only vehicles having a nonfunctional neutral (remove neutral for these)`r")
        break
    }
}

# 5. Replace the "Change neutral+clutch behavior in higher gears" paragraph
#    with three separate paragraphs.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*neutral+clutch*behavior in higher gears*") {
        $pRange = $p.Range
        $textRange = $d.Range($pRange.Start, $pRange.End - 1)
        $textRange.Text = "Cleaner vehicle swap/leave procedure"
        break
    }
}
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Change .ini format*") {
        $p.Range.InsertBefore("Change neutral/clutch revving to be more gradual and natural`r")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Change .ini format*") {
        $p.Range.InsertBefore("Change clutch slipping in higher gears to be less fake`r")
        break
    }
}

# 6. Collapse the two trailing empty paragraphs into one, and move the
#    _GoBack bookmark to that final paragraph.
$n = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($n - 1)
$secondToLast.Range.Delete()
$lastPara = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
